# edit.ps1 - applies the changes described by the diff:
#  1. Paragraph 1 ("This is a Microsoft word document.") gets two trailing
#     spaces appended to its existing run, followed by three new
#     red-colored (FF0000) runs spelling out:
#        "(This is a change " + EN DASH + " Ve" / "rsion for main branch" / ")"
#  2. A new, empty paragraph is appended at the very end of the document
#     (just before the sectPr) with paragraph shading fill F9F9F9
#     (<w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/>).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: modify the first paragraph.
# ---------------------------------------------------------------------------

$p1 = $d.Paragraphs(1)

# Append two trailing spaces to the existing (only) run's text, right before
# the paragraph mark.
$endOfP1Text = $p1.Range.End - 1
$spacesRange = $d.Range($endOfP1Text, $endOfP1Text)
$spacesRange.InsertAfter("  ")

# Build the en dash character used in "(This is a change <EN DASH> Ve".
$enDash = [char]0x2013

# Run 1: "(This is a change <EN DASH> Ve"  (red)
$e1 = $d.Paragraphs(1).Range.End - 1
$r1 = $d.Range($e1, $e1)
$r1.InsertAfter("(This is a change " + $enDash + " Ve")
$r1.Font.Color = 255

# Run 2: "rsion for main branch"  (red)
$e2 = $d.Paragraphs(1).Range.End - 1
$r2 = $d.Range($e2, $e2)
$r2.InsertAfter("rsion for main branch")
$r2.Font.Color = 255

# Run 3: ")"  (red)
$e3 = $d.Paragraphs(1).Range.End - 1
$r3 = $d.Range($e3, $e3)
$r3.InsertAfter(")")
$r3.Font.Color = 255

# ---------------------------------------------------------------------------
# Part 2: append a new, empty, shaded paragraph at the end of the document.
# ---------------------------------------------------------------------------
# A plain InsertParagraphAfter() at the end of the document copies the
# direct character/paragraph formatting of the last (content) paragraph onto
# the new paragraph mark, which we do not want. To get a clean paragraph we
# instead build one at the very start of the document (where there is
# nothing preceding it to inherit formatting from), move it to the end via
# the clipboard, and then apply only the shading we need.

# Create a clean paragraph mark at the very beginning of the document.
$startRange = $d.Range(0, 0)
$startRange.InsertParagraphAfter()

# Give it a temporary visible placeholder character so that Copy/Paste of
# the range actually carries content across (pasting a bare paragraph mark
# is a no-op in this environment).
$placeholderRange = $d.Range(0, 0)
$placeholderRange.InsertBefore("Z")

# Copy "Z" + the paragraph mark, then remove the temporary paragraph again.
$tempPara = $d.Paragraphs(1)
$copyRange = $d.Range($tempPara.Range.Start, $tempPara.Range.End)
$copyRange.Copy()
$copyRange.Delete()

# Paste at the very end of the document. The placeholder character merges
# onto the end of the last paragraph's text and the paragraph mark after it
# starts a brand new, cleanly-formatted paragraph.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.Paste()

# Remove the placeholder "Z" character that landed at the end of what is now
# the second-to-last paragraph.
$countAfterPaste = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs($countAfterPaste - 1)
$placeholderPos = $secondToLast.Range.End - 1
$cleanupRange = $d.Range($placeholderPos - 1, $placeholderPos)
$cleanupRange.Delete()

# Finally, reset the new last paragraph to the Normal style (clearing any
# remaining inherited paragraph formatting) and apply the desired shading.
$newLastPara = $d.Paragraphs($d.Paragraphs.Count)
$newLastPara.Range.Style = "Normal"
$newLastPara.Format.Shading.Texture = 0
$newLastPara.Format.Shading.ForegroundPatternColor = -16777216
$newLastPara.Format.Shading.BackgroundPatternColor = 0xF9F9F9
